# Insert a new row at position 345 on the (only) worksheet.
# This shifts the existing rows 345-430 down to 346-431, and the new
# row 345 is populated with a new data record (weekly Choclo price entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(345).Insert()

$ws.Range("A345").Value = 4
$ws.Range("B345").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C345").Value = "Los Lagos"
$ws.Range("D345").Value = 45244
$ws.Range("E345").Value = 10
$ws.Range("F345").Value = 100112024
$ws.Range("G345").Value = "Choclo"
$ws.Range("H345").Value = "Dulce o Americano"
$ws.Range("I345").Value = "Primera"
$ws.Range("J345").Value = 100
$ws.Range("K345").Value = 48000
$ws.Range("L345").Value = 48000
$ws.Range("M345").Value = 48000
$ws.Range("N345").Value = "`$/malla 70 unidades"
$ws.Range("O345").Value = "Región de Arica y Parinacota"
$ws.Range("P345").Value = 686
$ws.Range("Q345").Value = 70
$ws.Range("R345").Value = "Hortaliza"
